# The document's headers/footers each contain a single inline picture
# (the Pearson logo in both footers, the BTec logo in the second header).
# Their OOXML <wp:docPr>/<pic:cNvPr> "name" attributes were swapped:
#   - footer 1 and footer 2 pictures: image1.png -> image2.png
#   - header 2 picture:               image2.jpg -> image1.jpg
#
# Word's InlineShape object model exposes this as the shape's .Name
# property, so rename each inline picture accordingly.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer 1 - Pearson Edexcel logo: image1.png -> image2.png
$footer1Shape = $sec.Footers.Item(1).Range.InlineShapes.Item(1)
$footer1Shape.Name = "image2.png"

# Footer 2 - Pearson Edexcel logo: image1.png -> image2.png
$footer2Shape = $sec.Footers.Item(2).Range.InlineShapes.Item(1)
$footer2Shape.Name = "image2.png"

# Header 2 - BTec logo: image2.jpg -> image1.jpg
$header2Shape = $sec.Headers.Item(2).Range.InlineShapes.Item(1)
$header2Shape.Name = "image1.jpg"

Write-Output "Renamed footer1, footer2, and header2 inline picture shapes."
